$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> [$old]"
    }
    return $found
}

# 1. Heading3 title: "Genetics and Plant Biotechnology" -> "Plant Genetics and Biotechnology"
Replace-Text "Genetics and Plant Biotechnology" "Plant Genetics and Biotechnology"

# 2. Activation date
Replace-Text "Ativação: 01/01/2019" "Ativação: 01/01/2025"

# 3. Fill the empty italic "Objetivos" (English) paragraph with its text
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertBefore("Promoting understanding of plant biotechnology, encompassing plant biology, genetic manipulation techniques, and the practical applications of this technology in cultivating and growing genetically modified plants.")

# 4. Drop trailing period on the English "Programa resumido" sentence
Replace-Text "Science, society and environmental impact of GM crops." "Science, society and environmental impact of GM crops"

# 5. Portuguese "Programa" numbered list rewrite.
# (The target text ends with a literal straight double-quote character. Find/Replace's
#  replacement text goes through AutoCorrect's "smart quotes" and turns a typed `"` into a
#  curly U+201D, so the quote mark is appended afterwards via InsertAfter, which bypasses
#  AutoCorrect and keeps it a literal straight quote, U+0022.)
Replace-Text "1. Introdução em Biotecnologia Vegetal e Agricultura2. Fisiologia e Desenvolvimento Vegetal3. Cultura de Tecidos de Plantas4. Regulação Gênica em Eucariotos5. Transformação Genética de Plantas mediado por Agrobacterium6. Biolística e outros Métodos de Transformação Genética de Plantas7. Vetores de Transformação de Plantas8. Genes, Características de Interesse e Estratégias Bioctecnológicas para a Engenharia de Plantas9. Risco e Benefícios associados a Plantas Geneticamente Modificadas (GM)10. Biologia Sintética em Plantas" `
    "1. Introdução em Biotecnologia Vegetal e Agricultura2. Fisiologia e Desenvolvimento Vegetal3. Cultura de Tecidos de Plantas4. Regulação Gênica em Eucariotos5. Transformação Genética de Plantas mediado por Agrobacterium e Biobalística6. Vetores de Transformação Genética de Plantas7. Genes, Características de Interesse e Estratégias Bioctecnológicas para a Engenharia de Plantas8. Risco e Benefícios associados a Plantas Geneticamente Modificadas (GM)9. Viagem Didática Complementar"

$progPt = $d.Content.Duplicate
$progPt.Find.Execute("Viagem Didática Complementar")
$progPt.Collapse(0)
$progPt.InsertAfter([char]34)

# 6. English "Programa" numbered list rewrite
Replace-Text "1. Introduction in plant biotechnology and agriculture2. Plant physiology and development3. Plant tissue culture4. Eukaryotic gene regulation5. Agrobacterium-mediated plant genetic transformation6. Biolistic and other non-Agrobacterium technologies of plant transformation7. Vectors for plant transformation8. Genes, traits of interest and Biotechnological strategies for engineering plants9. Risk and benefits associated with genetically modified (GM) plants10. Synthetic biology in plants" `
    "1. Introduction in plant biotechnology and agriculture2. Plant physiology and development3. Plant tissue culture4. Eukaryotic gene regulation5. Plant transformation using Agrobacterium and biolistic6. Vectors for plant transformation7. Genes, traits of interest and Biotechnological strategies for engineering plants8. Risk and benefits associated with genetically modified (GM) plants9. Educational Excursion"

# 7. Requisitos list: insert a new LOT2040 run before the Bioquímica run (with new course text),
#    remove the old LOT2040 "Engenharia Genética" run entirely, and update the Microbiologia run's text.
$v = [char]11

$bioq = $d.Content.Duplicate
$bioq.Find.Execute("LOT2008 -  Bioquímica II  (Requisito fraco)")
$insPoint = $bioq.Duplicate
$insPoint.Collapse(1)
$insPoint.InsertBefore("LOT2040 -  Engenharia Genética Teórica e Prática  (Requisito fraco)" + $v)

$oldEng = $d.Content.Duplicate
$oldEng.Find.Execute("LOT2040 -  Engenharia Genética  (Requisito fraco)")
$oldEng.MoveEnd(1, 1)
$oldEng.Text = ""

Replace-Text "LOT2053 -  Microbiologia  (Requisito fraco)" "LOT2053 -  Microbiologia: da Teoria à Prática  (Requisito fraco)"
